$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#PROJECTNAME#")

# --- simple value corrections (planning vs actual hours) ---
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1

$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 6

$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 4

$ws.Range("H56").Value = 8

$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 5

# --- insert the new "Remote access" task before the Documentatie header row ---
$ws.Rows.Item(64).Insert()

$ws.Range("A64").Value = 39
$ws.Range("C64").Value = "Remote access "
$ws.Range("D64").Value = "Could"
$ws.Range("E64").Value = "Mboo"
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 2

# --- fill in the WBS items for the Documentatie section (rows 66-79) ---
$ws.Range("C66").Value = "Interview"
$ws.Range("D66").Value = "Must"
$ws.Range("E66").Value = "Mjak"
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 2

$ws.Range("C67").Value = "Project plan"
$ws.Range("D67").Value = "Must"
$ws.Range("E67").Value = "Mboo"
$ws.Range("G67").Value = 5
$ws.Range("H67").Value = 5

$ws.Range("C68").Value = "Kerntaken"
$ws.Range("D68").Value = "Must"
$ws.Range("E68").Value = "Mboo/Mjak"
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 6

$ws.Range("C69").Value = "Technisch ontwerp"
$ws.Range("D69").Value = "Must"
$ws.Range("E69").Value = "Mboo"
$ws.Range("G69").Value = 8
$ws.Range("H69").Value = 8
$ws.Range("I69").ClearContents()

$ws.Range("C70").Value = "Functioneel ontwerp"
$ws.Range("D70").Value = "Must"
$ws.Range("E70").Value = "Mboo"
$ws.Range("G70").Value = 8
$ws.Range("H70").Value = 8
$ws.Range("I70").ClearContents()

$ws.Range("C71").Value = "Planning"
$ws.Range("D71").Value = "Must"
$ws.Range("E71").Value = "Mjak"
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 4

$ws.Range("C72").Value = "Examen portfolio"
$ws.Range("D72").Value = "Must"
$ws.Range("E72").Value = "Mboo"
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 8

$ws.Range("C73").Value = "WBS "
$ws.Range("D73").Value = "Must"
$ws.Range("E73").Value = "Mboo"
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 2

$ws.Range("C74").Value = "Fancy poster"
$ws.Range("D74").Value = "Must"
$ws.Range("E74").Value = "Mjak"
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 4

$ws.Range("C75").Value = "Technisch poster"
$ws.Range("D75").Value = "Must"
$ws.Range("E75").Value = "Mboo"
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 4

$ws.Range("C76").Value = "CRPR lijst"
$ws.Range("D76").Value = "Must"
$ws.Range("E76").Value = "Mjak"
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 2

$ws.Range("C77").Value = "Test lijst"
$ws.Range("D77").Value = "Must"
$ws.Range("E77").Value = "Mboo"
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 2

$ws.Range("C78").Value = "Acceptatie test"
$ws.Range("D78").Value = "Must"
$ws.Range("E78").Value = "Mboo"
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 2

# --- insert 3 new rows before the subtotal row to host the last task + two blank rows ---
$ws.Range("A79:A81").EntireRow.Insert()

$ws.Range("A79").Value = 53
$ws.Range("C79").Value = "Evaluatie verslag"
$ws.Range("D79").Value = "Must"
$ws.Range("E79").Value = "Mjak"
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 4

$ws.Range("A80").Value = 54
$ws.Range("A81").Value = 55

# --- fix up the subtotal formulas now living on row 82 ---
$ws.Range("G82").Formula = "=SUBTOTAL(9,G7:G79)"
$ws.Range("H82").Formula = "=SUM(H7:H79)"

# --- workbook level metadata ---
$ws.AutoFilterMode = $false
$ws.Range("A4:J81").AutoFilter()

$wb.Names("_xlnm.Print_Area").RefersTo = "='#PROJECTNAME#'!`$A`$1:`$J`$82"
